$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 99.333336
$ws.Range("I4").Value = 99
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 99
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -328
$ws.Range("H21").Value = 6166.6665
$ws.Range("I21").Value = 6166.6665
$ws.Range("K21").Value = 6166.6665
$ws.Range("M21").Value = -5698.6665
$ws.Range("H23").Value = 6166.6665
$ws.Range("I23").Value = 6166.6665
$ws.Range("K23").Value = 6166.6665
$ws.Range("M23").Value = -5932.6665
$ws.Range("H62").Value = 3928.1428
$ws.Range("I62").Value = 3659.4
$ws.Range("K62").Value = 3659.4
$ws.Range("M62").Value = -3035.4
$ws.Range("H65").Value = 3928.1428
$ws.Range("I65").Value = 3659.4
$ws.Range("K65").Value = 18297
$ws.Range("M65").Value = -15177
$ws.Range("H107").Value = 213
$ws.Range("I107").Value = 213
$ws.Range("K107").Value = 213
$ws.Range("M107").Value = 1707
$ws.Range("H135").Value = 1031.1666
$ws.Range("J135").Value = 1100
$ws.Range("L135").Value = 9900
$ws.Range("N135").Value = -14970
$ws.Range("H141").Value = 1109.875
$ws.Range("I141").Value = 1109.875
$ws.Range("K141").Value = 3329.625
$ws.Range("M141").Value = 1850.375

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1099.5
$ws.Range("I102").Value = 900
$ws.Range("K102").Value = 900
$ws.Range("M102").Value = 722
$ws.Range("H130").Value = 79214
$ws.Range("J130").Value = 79214
$ws.Range("L130").Value = 79214
$ws.Range("N130").Value = -89254

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2320.889
$ws.Range("I134").Value = 2320.889
$ws.Range("K134").Value = 6962.667
$ws.Range("M134").Value = -4427.667

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 289.83334
$ws.Range("I5").Value = 220
$ws.Range("K5").Value = 220
$ws.Range("M5").Value = -108
$ws.Range("H11").Value = 730.7143
$ws.Range("I11").Value = 9
$ws.Range("J11").Value = 851
$ws.Range("K11").Value = 9
$ws.Range("L11").Value = 851
$ws.Range("M11").Value = 131
$ws.Range("N11").Value = -1131
$ws.Range("H12").Value = 601.25
$ws.Range("I12").Value = 405
$ws.Range("J12").Value = 666.6667
$ws.Range("K12").Value = 405
$ws.Range("L12").Value = 666.6667
$ws.Range("M12").Value = -235
$ws.Range("N12").Value = -1006.6667
$ws.Range("H13").Value = 782.5
$ws.Range("I13").Value = 350
$ws.Range("J13").Value = 926.6667
$ws.Range("K13").Value = 350
$ws.Range("L13").Value = 926.6667
$ws.Range("M13").Value = -211
$ws.Range("N13").Value = -1204.6667
$ws.Range("H64").Value = 30000
$ws.Range("I64").Value = 30000
$ws.Range("K64").Value = 30000
$ws.Range("M64").Value = -29752
$ws.Range("H67").Value = 30000
$ws.Range("I67").Value = 30000
$ws.Range("K67").Value = 30000
$ws.Range("M67").Value = -29142
$ws.Range("H132").Value = 7632.0835
$ws.Range("I132").Value = 4226.4287
$ws.Range("K132").Value = 12679.2861
$ws.Range("M132").Value = -10149.2861
$ws.Range("H141").Value = 421986.1
$ws.Range("J141").Value = 421986.1
$ws.Range("L141").Value = 421986.1
$ws.Range("N141").Value = -432346.1

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 4191.6665
$ws.Range("I11").Value = 1000
$ws.Range("K11").Value = 3000
$ws.Range("M11").Value = -2860
$ws.Range("H34").Value = 1099.8334
$ws.Range("I34").Value = 400
$ws.Range("J34").Value = 1239.8
$ws.Range("K34").Value = 1200
$ws.Range("L34").Value = 3719.4
$ws.Range("M34").Value = -1116
$ws.Range("N34").Value = -3887.4
$ws.Range("H39").Value = 3859.5715
$ws.Range("J39").Value = 3859.5715
$ws.Range("L39").Value = 11578.7145
$ws.Range("N39").Value = -12166.7145
$ws.Range("H55").Value = 3005
$ws.Range("J55").Value = 3005
$ws.Range("L55").Value = 9015
$ws.Range("N55").Value = -9369
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("H139").Value = 5286
$ws.Range("I139").Value = 5286
$ws.Range("K139").Value = 15858
$ws.Range("M139").Value = -10718
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 59
$ws.Range("I2").Value = 28.25
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 28.25
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 84.75
$ws.Range("N2").Value = -326
$ws.Range("H70").Value = 7503.75
$ws.Range("I70").Value = 5431.5713
$ws.Range("K70").Value = 5431.5713
$ws.Range("M70").Value = -5161.5713
$ws.Range("H73").Value = 7503.75
$ws.Range("I73").Value = 5431.5713
$ws.Range("K73").Value = 5431.5713
$ws.Range("M73").Value = -4495.5713
$ws.Range("H122").Value = 1199.2
$ws.Range("I122").Value = 1332.3334
$ws.Range("J122").Value = 999.5
$ws.Range("K122").Value = 3997.0002
$ws.Range("L122").Value = 2998.5
$ws.Range("M122").Value = -1547.0002
$ws.Range("N122").Value = -7898.5

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2795.8
$ws.Range("J2").Value = 989.5
$ws.Range("L2").Value = 989.5
$ws.Range("N2").Value = -1213.5
$ws.Range("H34").Value = 30000
$ws.Range("I34").Value = 30000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 30000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -29828
$ws.Range("N34").ClearContents()
$ws.Range("H122").Value = 3990
$ws.Range("I122").Value = 3990
$ws.Range("K122").Value = 11970
$ws.Range("M122").Value = -9520

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 738
$ws.Range("I132").Value = 650.6667
$ws.Range("K132").Value = 1952.0001
$ws.Range("M132").Value = 577.9999
